$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.876.63"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.51%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.904.83"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.49%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.009"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.75"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.66%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.007"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.31%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4839"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.58%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3795"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07369"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.77%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9307"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.13%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.73"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.33%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.967.25"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.84%  "
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07743"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.59%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.487"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.70%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.625"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.01%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.84"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.90%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.008"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.33%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008850"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.39%  "
$ws.Range("E19").Value = "  -0.36%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "27.948.47"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.81%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.67"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.25%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.156"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.30%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.182.98"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.61%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.89"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.85%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "156.04"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.70%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.918"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.10%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.45"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.01%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.132"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +6.32%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "117.29"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.73%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.952"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.17%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08959"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.96%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.205"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.72%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.255"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.23%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7645"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.83%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.653"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.70%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02049"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.20%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.527"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.35%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.095"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.34%  "
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.005"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.83%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5475"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.39%  "
$ws.Range("B41").Value = "Hedera"
$ws.Range("C41").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.05268"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.24%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.951"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.75%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1527"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.69%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.459"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.12%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.70"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.43%  "
$ws.Range("B46").Value = "Quant"
$ws.Range("C46").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "109.20"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.26%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4809"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.16%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.006"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.43%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.651"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.74%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "67.95"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.06%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06090"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.08%  "
